$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 into the new I1:J1 header cells, then set their text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I40 with 1, and J2:J40 with the same value as the corresponding H cell
for ($r = 2; $r -le 40; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
